# Applies the "Internal Keycode" / "Binding Keycode" columns (E, F) to the
# Keybindings Guide worksheet, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths for the two new columns (E, F).
# Excel's ColumnWidth setter quantizes to 1/6-character increments, so these
# are the closest settable values that reproduce the target stored widths
# (16 and ~15.285) in the saved XML.
# ---------------------------------------------------------------------------
$ws.Range("E1").ColumnWidth = 15.1666666666667
$ws.Range("F1").ColumnWidth = 14.5

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("E1").Value = "Internal Keycode"
$ws.Range("F1").Value = "Binding Keycode"
$ws.Range("E1:F1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Data rows: Internal Keycode (E) and Binding Keycode (F)
# Most Binding Keycode cells simply mirror the Internal Keycode number, but a
# few rows instead carry a descriptive, right-aligned text value.
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 13

# Row 3 - Binding Keycode is "N/A" (italic, right aligned)
$ws.Range("E3").Value = 27
$ws.Range("F3").Value = "N/A"
$ws.Range("F3").HorizontalAlignment = -4152
$ws.Range("F3").Font.Italic = $true

# Row 4
$ws.Range("E4").Value = 36
$ws.Range("F4").Value = 27

# Row 5
$ws.Range("E5").Value = 38
$ws.Range("F5").Value = 38

# Row 6
$ws.Range("E6").Value = 40
$ws.Range("F6").Value = 40

# Row 7
$ws.Range("E7").Value = 37
$ws.Range("F7").Value = 37

# Row 8
$ws.Range("E8").Value = 39
$ws.Range("F8").Value = 39

# Row 9
$ws.Range("E9").Value = 87
$ws.Range("F9").Value = 87

# Row 10
$ws.Range("E10").Value = 83
$ws.Range("F10").Value = 83

# Row 11
$ws.Range("E11").Value = 65
$ws.Range("F11").Value = 65

# Row 12
$ws.Range("E12").Value = 68
$ws.Range("F12").Value = 68

# Row 13 - Binding Keycode is a descriptive text, right aligned (normal font)
$ws.Range("E13").Value = 81
$ws.Range("F13").Value = "81, 32"
$ws.Range("F13").HorizontalAlignment = -4152

# Row 14 - Binding Keycode keeps numeric value, but right aligned
$ws.Range("E14").Value = 69
$ws.Range("F14").Value = 69
$ws.Range("F14").HorizontalAlignment = -4152

# Row 15
$ws.Range("E15").Value = 56
$ws.Range("F15").Value = "56, 104"
$ws.Range("F15").HorizontalAlignment = -4152

# Row 16 - Internal Keybinding text (D) gains extra binding info
$ws.Range("D16").Value = "5, numpad 5, 2, numpad 2"
$ws.Range("E16").Value = 53
$ws.Range("F16").Value = "53, 101, 50, 98"
$ws.Range("F16").HorizontalAlignment = -4152

# Row 17
$ws.Range("D17").Value = "4, numpad 4"
$ws.Range("E17").Value = 52
$ws.Range("F17").Value = "52, 100"
$ws.Range("F17").HorizontalAlignment = -4152

# Row 18
$ws.Range("D18").Value = "6, numpad 6"
$ws.Range("E18").Value = 54
$ws.Range("F18").Value = "54, 102"
$ws.Range("F18").HorizontalAlignment = -4152

# Row 19
$ws.Range("E19").Value = 73
$ws.Range("F19").Value = 73

# Row 20
$ws.Range("E20").Value = 75
$ws.Range("F20").Value = 75

# Row 21
$ws.Range("E21").Value = 74
$ws.Range("F21").Value = 74

# Row 22
$ws.Range("E22").Value = 76
$ws.Range("F22").Value = 76

# Row 23
$ws.Range("E23").Value = 85
$ws.Range("F23").Value = 85

# Row 24
$ws.Range("E24").Value = 79
$ws.Range("F24").Value = 79

# ---------------------------------------------------------------------------
# Selection, matching the author's final cursor position in the diff.
# ---------------------------------------------------------------------------
$ws.Range("E23").Select()
